$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update refreshed crypto "Price" (D) and "Volume(1h)" (E) values.
# Values are prefixed with a leading apostrophe so Excel stores them as
# text (matching the source workbook's string cells) instead of
# auto-converting numeric-looking strings (e.g. "597.74") into floats.
# ClearFormats() afterwards drops the transient quote-prefix style so the
# cells keep their original (unstyled) appearance.
$ws.Range("D2").Value = '''68.341.43'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '''  +0.17%  '
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = '''2.651.36'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '''  +0.54%  '
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = '''  -0.03%  '
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = '''597.74'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '''  -0.20%  '
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = '''158.80'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '''  +2.92%  '
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = '''  -0.04%  '
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = '''  -0.31%  '
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = '''  +3.71%  '
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = '''  -1.19%  '
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = '''  +0.55%  '
$ws.Range("E11").ClearFormats()
$ws.Range("E12").Value = '''  +0.83%  '
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = '''28.12'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '''  +0.56%  '
$ws.Range("E13").ClearFormats()
$ws.Range("E14").Value = '''  +1.37%  '
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = '''3.130.76'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '''  +0.49%  '
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = '''68.222.95'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '''  +0.02%  '
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = '''2.673.01'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '''  +0.96%  '
$ws.Range("E17").ClearFormats()
$ws.Range("E18").Value = '''  +0.07%  '
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = '''364.55'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '''  -0.38%  '
$ws.Range("E19").ClearFormats()
$ws.Range("E20").Value = '''  -0.74%  '
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = '''  +3.66%  '
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = '''4.83'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '''  -0.35%  '
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = '''  -1.88%  '
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = '''75.21'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '''  +2.18%  '
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = '''  +0.24%  '
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = '''9.75'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '''  -2.85%  '
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = '''2.784.33'
$ws.Range("D27").ClearFormats()
$ws.Range("E28").Value = '''  +0.46%  '
$ws.Range("E28").ClearFormats()
$ws.Range("E29").Value = '''  -0.01%  '
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = '''558.92'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '''  -2.58%  '
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = '''  +0.42%  '
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = '''  -0.43%  '
$ws.Range("E32").ClearFormats()
$ws.Range("E33").Value = '''  +0.46%  '
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = '''  -0.74%  '
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = '''0.999'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '''  -0.03%  '
$ws.Range("E35").ClearFormats()
$ws.Range("E36").Value = '''  +1.94%  '
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = '''19.83'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '''  +2.98%  '
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = '''159.67'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '''  -0.35%  '
$ws.Range("E38").ClearFormats()
$ws.Range("E39").Value = '''  +1.14%  '
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = '''  -2.19%  '
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = '''5.36'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '''  -0.29%  '
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = '''0.0₆0332'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '''  +3.35%  '
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = '''2.63'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '''  +0.02%  '
$ws.Range("E43").ClearFormats()
$ws.Range("E44").Value = '''  +0.05%  '
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = '''158.52'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '''  +1.03%  '
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = '''  +0.01%  '
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = '''  +1.44%  '
$ws.Range("E47").ClearFormats()
$ws.Range("E48").Value = '''  -0.80%  '
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = '''  -0.07%  '
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = '''0.616'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '''  +0.25%  '
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = '''  +1.11%  '
$ws.Range("E51").ClearFormats()
